# Update cryptocurrency price/volume data (and one rebranded coin row)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '60.121.11'
$ws.Range('E2').Value = '  -0.36%  '
$ws.Range('D3').Value = '2.418.91'
$ws.Range('E3').Value = '  -0.53%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range("D5").Value = "'552.39"
$ws.Range('E5').Value = '  -0.67%  '
$ws.Range("D6").Value = "'137.16"
$ws.Range('E6').Value = '  -1.66%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range("D8").Value = "'0.589"
$ws.Range('E8').Value = '  +1.87%  '
$ws.Range("D9").Value = "'0.105"
$ws.Range('E9').Value = '  -1.69%  '
$ws.Range("D10").Value = "'5.65"
$ws.Range('E10').Value = '  -1.83%  '
$ws.Range('E11').Value = '  -0.32%  '
$ws.Range("D12").Value = "'0.353"
$ws.Range('E12').Value = '  -1.57%  '
$ws.Range("D13").Value = "'24.91"
$ws.Range('E13').Value = '  -0.49%  '
$ws.Range('D14').Value = '2.853.47'
$ws.Range('E14').Value = '  -0.37%  '
$ws.Range('D15').Value = '60.044.43'
$ws.Range('E15').Value = '  -0.34%  '
$ws.Range('E16').Value = '  -1.57%  '
$ws.Range('D17').Value = '2.420.82'
$ws.Range('E17').Value = '  -1.12%  '
$ws.Range("D18").Value = "'11.30"
$ws.Range('E18').Value = '  -0.82%  '
$ws.Range("D19").Value = "'4.48"
$ws.Range('E19').Value = '  +1.00%  '
$ws.Range("D20").Value = "'328.03"
$ws.Range('E20').Value = '  -1.85%  '
$ws.Range("D21").Value = "'6.74"
$ws.Range('E21').Value = '  -0.45%  '
$ws.Range("D22").Value = "'0.996"
$ws.Range('E22').Value = '  -0.33%  '
$ws.Range("D23").Value = "'65.53"
$ws.Range('E23').Value = '  +0.39%  '
$ws.Range("D24").Value = "'0.177"
$ws.Range('E24').Value = '  +2.87%  '
$ws.Range("D25").Value = "'8.70"
$ws.Range('E25').Value = '  +0.63%  '
$ws.Range("D26").Value = "'1.02"
$ws.Range('E26').Value = '  +1.57%  '
$ws.Range("D27").Value = "'1.39"
$ws.Range('E27').Value = '  +2.20%  '
$ws.Range('D28').Value = '0.0₃0773'
$ws.Range('E28').Value = '  -2.36%  '
$ws.Range("D29").Value = "'1.75"
$ws.Range('E29').Value = '  -2.35%  '
$ws.Range("D30").Value = "'169.94"
$ws.Range("D31").Value = "'6.10"
$ws.Range('E31').Value = '  -3.92%  '
$ws.Range('E32').Value = '  +2.62%  '
$ws.Range("D33").Value = "'0.404"
$ws.Range('E33').Value = '  -3.96%  '
$ws.Range("D34").Value = "'18.57"
$ws.Range('E34').Value = '  -1.14%  '
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range("D36").Value = "'1.32"
$ws.Range('E36').Value = '  +0.71%  '
$ws.Range("D38").Value = "'4.19"
$ws.Range('E38').Value = '  -1.15%  '
$ws.Range("D39").Value = "'330.17"
$ws.Range('E39').Value = '  +1.85%  '
$ws.Range("D40").Value = "'1.60"
$ws.Range('E40').Value = '  -1.30%  '
$ws.Range("D41").Value = "'38.87"
$ws.Range('E41').Value = '  -2.39%  '
$ws.Range("D42").Value = "'144.56"
$ws.Range('E42').Value = '  +2.60%  '
$ws.Range("D43").Value = "'3.65"
$ws.Range('E43').Value = '  -1.88%  '
$ws.Range("D44").Value = "'20.08"
$ws.Range('E44').Value = '  +1.98%  '
$ws.Range("D45").Value = "'0.0966"
$ws.Range('E45').Value = '  +0.40%  '
$ws.Range("D46").Value = "'0.0515"
$ws.Range('E46').Value = '  -2.30%  '
$ws.Range("D47").Value = "'0.577"
$ws.Range('E47').Value = '  +0.59%  '
$ws.Range("D48").Value = "'0.0223"
$ws.Range('E48').Value = '  -1.96%  '
$ws.Range("D49").Value = "'11.05"
$ws.Range('E49').Value = '  -0.02%  '
$ws.Range('E50').Value = '  -3.62%  '
$ws.Range('B51').Value = 'BitgetToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb'
$ws.Range("D51").Value = "'0.947"
$ws.Range('E51').Value = '  -0.50%  '
